$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.05619466666666667
$ws.Range("I2").Value = 0.04986276087265156
$ws.Range("J2").Value = 0.07297477932340853
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.016657
$ws.Range("N2").Value = 0.033314
$ws.Range("Q2").Value = 0.0009360345626666668
$ws.Range("R2").Value = 0.005616207376000001
$ws.Range("S2").Value = 0.04986276087265156
$ws.Range("T2").Value = 0.07297477932340853

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.070792
$ws.Range("H3").Value = 2.141584
$ws.Range("I3").Value = 0.9501372391273485
$ws.Range("J3").Value = 0.9270252206765914
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.016657
$ws.Range("N3").Value = 0.033314
$ws.Range("Q3").Value = 0.017836182344
$ws.Range("R3").Value = 0.071344729376
$ws.Range("S3").Value = 0.9501372391273485
$ws.Range("T3").Value = 0.9270252206765914
